# Add "preparer" (libraryPreparer) and "purpose" information to the data
# rows on Sheet1. Column B = libraryPreparer, Column E = purpose.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 29 }

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = "H.BROWN"
    $ws.Cells.Item($r, 5).Value = "fullRNASEQ"
}

# Select the range near the bottom of the data, matching the saved view state.
$ws.Range("A30:D33").Select()
